# Improved wishlistApp() and added rowCounter(String sheetName).
#
# rowCounter(sheetName) returns how many rows of the named sheet are
# currently in use, so wishlistApp() knows where the existing wishlist
# data ends.
#
# wishlistApp() rebuilds the Product sheet's wishlist rows (column A is the
# item id, column B is the item name). The old rows are cleared first (so
# rewritten cells fall back to the sheet's default formatting) and then
# rewritten: the two PAX entries move to the top of the list, a new
# LÅNGFJÄLL item is added, and the LINNMON/ALEX entry that used to be first
# is pushed down to the new last row.

function rowCounter([string]$sheetName) {
    $wb = $excel.ActiveWorkbook
    $sheet = $wb.Worksheets.Item($sheetName)
    return $sheet.UsedRange.Rows.Count
}

function wishlistApp() {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.Worksheets.Item("Product")

    $firstDataRow = 3
    $lastDataRow = rowCounter("Product")

    # Wipe the existing wishlist rows (values + formatting) before
    # rewriting them.
    $ws.Range("A" + $firstDataRow + ":B" + $lastDataRow).Clear()

    $wishlist = @(
        @(19288134, "PAX"),
        @(59931463, "PAX"),
        @(39252515, "LÅNGFJÄLL"),
        @(69222616, "LINNMON/ALEX")
    )

    $row = $firstDataRow
    foreach ($entry in $wishlist) {
        $ws.Cells.Item($row, 1).Value = $entry[0]
        $ws.Cells.Item($row, 2).Value = $entry[1]
        $row = $row + 1
    }

    $ws.Cells.Item($row - 1, 2).Select() | Out-Null
}

wishlistApp
